$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1500
$ws1.Range("F5").Value = 148
$ws1.Range("F9").Value = 138
$ws1.Range("F10").Value = 728
$ws1.Range("F13").Value = 319
$ws1.Range("F15").Value = 6339
$ws1.Range("F20").Value = 15191
$ws1.Range("F22").Value = 273
$ws1.Range("F23").Value = 135
$ws1.Range("F25").Value = 11004
$ws1.Range("F26").Value = 738
$ws1.Range("F27").Value = 4294
$ws1.Range("F28").Value = 225
$ws1.Range("F30").Value = 12
$ws1.Range("F31").Value = 125

# Sheet "全部类型" (fourth sheet) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1500
$ws4.Range("F5").Value = 148
$ws4.Range("F10").Value = 138
$ws4.Range("F11").Value = 728
$ws4.Range("F15").Value = 319
$ws4.Range("F18").Value = 6339
$ws4.Range("F23").Value = 15191
$ws4.Range("F25").Value = 273
$ws4.Range("F26").Value = 135
$ws4.Range("F28").Value = 11004
$ws4.Range("F29").Value = 738
$ws4.Range("F30").Value = 4294
$ws4.Range("F31").Value = 225
$ws4.Range("F33").Value = 12
$ws4.Range("F34").Value = 125
